$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to stay text (values like "305.00" would otherwise be
# auto-converted to numbers by Excel, losing the trailing zero / matching the
# original inline-string "Price" column formatting).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '41.942.19'
$ws.Range("E2").Value = '  +0.33%  '
$ws.Range("D3").Value = '2.260.40'
$ws.Range("E3").Value = '  -0.38%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '305.00'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").Value = '95.25'
$ws.Range("E6").Value = '  +2.51%  '
$ws.Range("E7").Value = '  -0.62%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +0.46%  '
$ws.Range("D10").Value = '35.02'
$ws.Range("E10").Value = '  +7.17%  '
$ws.Range("E11").Value = '  -1.15%  '
$ws.Range("E12").Value = '  -0.47%  '
$ws.Range("D13").Value = '6.62'
$ws.Range("E13").Value = '  -0.92%  '
$ws.Range("D14").Value = '2.612.24'
$ws.Range("E14").Value = '  -0.27%  '
$ws.Range("D15").Value = '14.39'
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").Value = '2.252.28'
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").Value = '41.857.48'
$ws.Range("E18").Value = '  +0.26%  '
$ws.Range("D19").Value = '12.34'
$ws.Range("E19").Value = '  -4.58%  '
$ws.Range("D20").Value = '0.0₃0900'
$ws.Range("E20").Value = '  -1.55%  '
$ws.Range("D21").Value = '5.95'
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("D22").Value = '67.64'
$ws.Range("E22").Value = '  -0.56%  '
$ws.Range("D23").Value = '237.13'
$ws.Range("E23").Value = '  -2.89%  '
$ws.Range("D24").Value = '2.57'
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  +0.02%  '
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").Value = '23.62'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '36.49'
$ws.Range("E28").Value = '  +5.14%  '
$ws.Range("D29").Value = '9.48'
$ws.Range("E29").Value = '  -1.77%  '
$ws.Range("E30").Value = '  +1.65%  '
$ws.Range("D31").Value = '160.08'
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("D32").Value = '5.20'
$ws.Range("E32").Value = '  -2.69%  '
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("D34").Value = '3.18'
$ws.Range("E34").Value = '  +4.91%  '
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").Value = '17.05'
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("E38").Value = '  -0.66%  '
$ws.Range("E39").Value = '  +1.46%  '
$ws.Range("E40").Value = '  -2.36%  '
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("D42").Value = '2.40'
$ws.Range("E42").Value = '  +6.72%  '
$ws.Range("D43").Value = '1.975.20'
$ws.Range("E43").Value = '  -1.80%  '
$ws.Range("D44").Value = '0.0282'
$ws.Range("E44").Value = '  +0.01%  '
$ws.Range("D45").Value = '18.83'
$ws.Range("E45").Value = '  -4.63%  '
$ws.Range("D46").Value = '2.93'
$ws.Range("E46").Value = '  +0.43%  '
$ws.Range("D47").Value = '9.88'
$ws.Range("E47").Value = '  -3.77%  '
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("D49").Value = '72.18'
$ws.Range("E49").Value = '  -1.40%  '
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '90.69'
$ws.Range("E51").Value = '  -1.30%  '

# Restore the default (unstyled) cell style now that the text values are
# committed, matching the original workbook where these cells carry no
# explicit style index.
$ws.Range("D2:D51").Style = "Normal"
